$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.219.81"
$ws.Range("E2").Value = "  +3.99%  "

$ws.Range("D3").Value = "'2.427.94"
$ws.Range("E3").Value = "  +3.28%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'554.04"
$ws.Range("E5").Value = "  +2.35%  "

$ws.Range("D6").Value = "'139.22"
$ws.Range("E6").Value = "  +3.39%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'0.578"
$ws.Range("E8").Value = "  +1.22%  "

$ws.Range("E9").Value = "  +3.54%  "

$ws.Range("D10").Value = "'5.77"

$ws.Range("D11").Value = "'0.358"
$ws.Range("E11").Value = "  +0.45%  "

$ws.Range("E12").Value = "  -2.12%  "

$ws.Range("D13").Value = "'25.02"
$ws.Range("E13").Value = "  +5.40%  "

$ws.Range("D14").Value = "'2.858.06"
$ws.Range("E14").Value = "  +3.25%  "

$ws.Range("D15").Value = "'60.119.46"
$ws.Range("E15").Value = "  +3.95%  "

$ws.Range("E16").Value = "  +3.87%  "

$ws.Range("D17").Value = "'2.428.80"
$ws.Range("E17").Value = "  +3.32%  "

$ws.Range("D18").Value = "'11.40"
$ws.Range("E18").Value = "  +6.37%  "

$ws.Range("D19").Value = "'4.41"
$ws.Range("E19").Value = "  +2.85%  "

$ws.Range("D20").Value = "'333.39"
$ws.Range("E20").Value = "  +1.03%  "

$ws.Range("E21").Value = "  +1.13%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").Value = "'65.23"
$ws.Range("E23").Value = "  +4.26%  "

$ws.Range("D24").Value = "'0.171"
$ws.Range("E24").Value = "  +3.12%  "

$ws.Range("D25").Value = "'8.63"
$ws.Range("E25").Value = "  +2.69%  "

$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.53%  "

$ws.Range("E27").Value = "  -0.43%  "

$ws.Range("D28").Value = "'0.0₃0793"
$ws.Range("E28").Value = "  +7.68%  "

$ws.Range("D29").Value = "'1.78"
$ws.Range("E29").Value = "  +1.35%  "

$ws.Range("D30").Value = "'6.33"
$ws.Range("E30").Value = "  +3.35%  "

$ws.Range("D31").Value = "'169.77"
$ws.Range("E31").Value = "  -0.41%  "

$ws.Range("E32").Value = "  +2.79%  "

$ws.Range("D33").Value = "'18.74"
$ws.Range("E33").Value = "  +2.11%  "

$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("E35").Value = "  +6.01%  "

$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("E37").Value = "  +0.15%  "

$ws.Range("E38").Value = "  +0.87%  "

$ws.Range("D39").Value = "'323.99"
$ws.Range("E39").Value = "  +11.90%  "

$ws.Range("D40").Value = "'0.420"
$ws.Range("E40").Value = "  +11.31%  "

$ws.Range("D41").Value = "'39.55"
$ws.Range("E41").Value = "  +1.31%  "

$ws.Range("E42").Value = "  +1.58%  "

$ws.Range("D43").Value = "'140.16"
$ws.Range("E43").Value = "  -1.71%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "'0.0962"
$ws.Range("E44").Value = "  +1.27%  "

$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0525"
$ws.Range("E45").Value = "  +3.21%  "

$ws.Range("D46").Value = "'19.59"
$ws.Range("E46").Value = "  +1.96%  "

$ws.Range("E47").Value = "  +8.34%  "

$ws.Range("E48").Value = "  +1.51%  "

$ws.Range("E49").Value = "  +2.00%  "

$ws.Range("D50").Value = "'17.87"
$ws.Range("E50").Value = "  +2.24%  "

$ws.Range("E51").Value = "  -0.12%  "
